$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6, shifting existing row 6 (hostname/adeye03u) down to row 7
$ws.Rows.Item(6).Insert()

$ws.Range("B6").Value = "/opt/ros/kinetic"
$ws.Range("A6").Value = "ROS_folder"

$ws.Range("A6").Select()
